$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 4) for PR-003
$ws.Range("A4").Value = 43343
$ws.Range("A4").NumberFormat = "d-mmm-yy"
$ws.Range("D4").Value = "PR-003"
$ws.Range("E4").Value = "Open"
$ws.Range("F4").Value = "Software-app"

# Update selection to match the authored state (F5 selected)
$ws.Range("F5").Select()
